# Auto-generated crypto price/volume update script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.544.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'1.876.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'0.7222"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.63%  "
$ws.Range("D6").Value = "'239.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.90%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -3.96%  "
$ws.Range("D9").Value = "'0.3086"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.63%  "
$ws.Range("D10").Value = "'25.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.71%  "
$ws.Range("D11").Value = "'0.08243"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").Value = "'1.904.03"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.98%  "
$ws.Range("D13").Value = "'0.7263"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.72%  "
$ws.Range("D14").Value = "'5.252"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").Value = "'90.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").Value = "'29.617.69"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.26%  "
$ws.Range("D17").Value = "'5.847"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.06%  "
$ws.Range("D18").Value = "'242.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.86%  "
$ws.Range("D19").Value = "'0.000007858"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.38%  "
$ws.Range("D20").Value = "'13.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").Value = "'2.150.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'7.806"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.46%  "
$ws.Range("D25").Value = "'0.1578"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.66%  "
$ws.Range("D26").Value = "'162.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("D27").Value = "'8.997"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("D28").Value = "'18.36"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Value = "'1.946"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("D30").Value = "'1.355"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.01%  "
$ws.Range("D31").Value = "'1.482"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").Value = "'4.346"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.87%  "
$ws.Range("D33").Value = "'4.078"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("D34").Value = "'0.05256"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("D35").Value = "'1.200"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.65%  "
$ws.Range("D36").Value = "'0.7190"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.59%  "
$ws.Range("D37").Value = "'0.9999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").Value = "'2.672"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").Value = "'0.01871"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.19%  "
$ws.Range("D40").Value = "'2.717"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("D41").Value = "'1.181.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.45%  "
$ws.Range("D42").Value = "'0.9118"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("D43").Value = "'6.010"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.38%  "
$ws.Range("E44").Value = "  +1.33%  "
$ws.Range("D45").Value = "'71.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.17%  "
$ws.Range("E46").Value = "  +0.12%  "
$ws.Range("D47").Value = "'103.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.33%  "
$ws.Range("D48").Value = "'0.5354"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.28%  "
$ws.Range("D49").Value = "'1.781"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.52%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.237"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").Value = "'7.076"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.90%  "
